$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.682.79'
$ws.Range('E2').Value = '  -0.57%  '
$ws.Range('D3').Value = '2.530.18'
$ws.Range('E3').Value = '  +1.32%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'544.91"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('D6').Value = "'146.20"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.80%  '
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').Value = "'0.575"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.37%  '
$ws.Range('D9').Value = '2.554.06'
$ws.Range('E9').Value = '  +1.15%  '
$ws.Range('E10').Value = '  +0.09%  '
$ws.Range('E11').Value = '  +0.19%  '
$ws.Range('D12').Value = "'5.60"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.54%  '
$ws.Range('D13').Value = "'0.361"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.15%  '
$ws.Range('D14').Value = '2.973.52'
$ws.Range('E14').Value = '  +0.88%  '
$ws.Range('D15').Value = "'23.65"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.18%  '
$ws.Range('D16').Value = '59.600.37'
$ws.Range('E16').Value = '  -0.74%  '
$ws.Range('D17').Value = "'0.0000143"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.50%  '
$ws.Range('D18').Value = '2.543.23'
$ws.Range('E18').Value = '  +1.17%  '
$ws.Range('D19').Value = "'11.27"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.83%  '
$ws.Range('D20').Value = "'4.31"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.69%  '
$ws.Range('D21').Value = "'327.54"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('D22').Value = "'0.999"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').Value = "'5.95"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.99%  '
$ws.Range('D24').Value = "'62.35"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.01%  '
$ws.Range('D25').Value = "'0.439"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.33%  '
$ws.Range('D26').Value = "'0.166"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.85%  '
$ws.Range('D27').Value = "'0.993"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.86%  '
$ws.Range('D28').Value = "'8.05"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.28%  '
$ws.Range('D29').Value = '0.0₃0802'
$ws.Range('E29').Value = '  +0.51%  '
$ws.Range('D30').Value = "'6.87"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.64%  '
$ws.Range('D31').Value = "'1.83"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.25%  '
$ws.Range('E32').Value = '  -7.71%  '
$ws.Range('E33').Value = '  +4.80%  '
$ws.Range('D34').Value = "'160.73"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.13%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').Value = "'18.79"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.94%  '
$ws.Range('D37').Value = "'4.43"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.49%  '
$ws.Range('D38').Value = "'1.63"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.72%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = "'5.71"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.54%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = "'37.13"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.88%  '
$ws.Range('D41').Value = "'0.850"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.09%  '
$ws.Range('D42').Value = "'299.05"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.18%  '
$ws.Range('D43').Value = "'3.72"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.14%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value = "'0.609"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.51%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = "'0.992"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.38%  '
$ws.Range('D46').Value = "'10.79"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.09%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'19.02"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.44%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = "'0.0940"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.79%  '
$ws.Range('D49').Value = "'123.78"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.90%  '
$ws.Range('D50').Value = "'0.0518"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.20%  '
$ws.Range('D51').Value = "'0.0229"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.75%  '
